# Price/volume-label/hour refresh for the crypto symbol list (GitHub Actions run).
# Every numeric-looking value is set with a leading apostrophe so Excel keeps it
# as text (matching the original inlineStr cells) instead of coercing to a number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'245.91"
$ws.Range("G2").Value = "'17"

# Row 3
$ws.Range("D3").Value = "'21.53"
$ws.Range("G3").Value = "'17"

# Row 4
$ws.Range("D4").Value = "'5.328"
$ws.Range("G4").Value = "'17"

# Row 5
$ws.Range("D5").Value = "'0.05656"
$ws.Range("G5").Value = "'17"

# Row 6
$ws.Range("D6").Value = "'3.384"
$ws.Range("G6").Value = "'17"

# Row 7
$ws.Range("D7").Value = "'0.8063"
$ws.Range("G7").Value = "'17"

# Row 8
$ws.Range("G8").Value = "'17"

# Row 9
$ws.Range("D9").Value = "'0.1457"
$ws.Range("G9").Value = "'17"

# Row 10
$ws.Range("D10").Value = "'0.07431"
$ws.Range("G10").Value = "'17"

# Row 11
$ws.Range("D11").Value = "'0.03172"
$ws.Range("G11").Value = "'17"

# Row 12
$ws.Range("D12").Value = "'0.03105"
$ws.Range("G12").Value = "'17"

# Row 13
$ws.Range("D13").Value = "'0.09263"
$ws.Range("G13").Value = "'17"

# Row 14
$ws.Range("D14").Value = "'3.568"
$ws.Range("G14").Value = "'17"

# Row 15
$ws.Range("D15").Value = "'0.001658"
$ws.Range("G15").Value = "'17"

# Row 16
$ws.Range("D16").Value = "'0.04719"
$ws.Range("G16").Value = "'17"

# Row 17
$ws.Range("D17").Value = "'0.0005826"
$ws.Range("E17").Value = "'16OneONE"
$ws.Range("G17").Value = "'17"

# Row 18
$ws.Range("D18").Value = "'0.006372"
$ws.Range("G18").Value = "'17"

# Row 19
$ws.Range("D19").Value = "'0.005001"
$ws.Range("G19").Value = "'17"

# Row 20
$ws.Range("G20").Value = "'17"

# Row 21
$ws.Range("D21").Value = "'0.0001501"
$ws.Range("G21").Value = "'17"

# Row 22
$ws.Range("D22").Value = "'0.0003103"
$ws.Range("G22").Value = "'17"

# Row 23
$ws.Range("D23").Value = "'3.770"
$ws.Range("G23").Value = "'17"

# Row 24
$ws.Range("D24").Value = "'6.394"
$ws.Range("G24").Value = "'17"

# Row 25
$ws.Range("D25").Value = "'2.097"
$ws.Range("G25").Value = "'17"

# Row 26
$ws.Range("D26").Value = "'0.3288"
$ws.Range("G26").Value = "'17"

# Row 27
$ws.Range("G27").Value = "'17"

# Row 28
$ws.Range("G28").Value = "'17"

# Row 29
$ws.Range("G29").Value = "'17"

# Row 30
$ws.Range("G30").Value = "'17"

# Row 31
$ws.Range("G31").Value = "'17"

# Row 32
$ws.Range("G32").Value = "'17"

# Row 33
$ws.Range("G33").Value = "'17"

# Row 34
$ws.Range("G34").Value = "'17"

# Row 35
$ws.Range("G35").Value = "'17"

# Row 36
$ws.Range("G36").Value = "'17"

# Row 37
$ws.Range("G37").Value = "'17"

# Row 38
$ws.Range("G38").Value = "'17"

# Row 39
$ws.Range("G39").Value = "'17"

# Row 40
$ws.Range("D40").Value = "'0.04009"
$ws.Range("G40").Value = "'17"

# Row 41
$ws.Range("D41").Value = "'0.006986"
$ws.Range("E41").Value = "'40KickTokenKICKBestin24h"
$ws.Range("G41").Value = "'17"

# Row 42
$ws.Range("D42").Value = "'0.003503"
$ws.Range("G42").Value = "'17"

# Row 43
$ws.Range("G43").Value = "'17"

# Row 44
$ws.Range("D44").Value = "'0.008473"
$ws.Range("G44").Value = "'17"

# Row 45
$ws.Range("D45").Value = "'0.00005948"
$ws.Range("G45").Value = "'17"

# Row 46
$ws.Range("G46").Value = "'17"

# Row 47
$ws.Range("D47").Value = "'0.0005505"
$ws.Range("E47").Value = "'46ACDXExchangeACXTWorstin24h"
$ws.Range("G47").Value = "'17"

# Row 48
$ws.Range("D48").Value = "'0.6831"
$ws.Range("G48").Value = "'17"

# Row 49
$ws.Range("D49").Value = "'0.02640"
$ws.Range("E49").Value = "'48BOLOBOLO"
$ws.Range("G49").Value = "'17"

# Row 50
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("G50").Value = "'17"

# Row 51
$ws.Range("G51").Value = "'17"
